# Mise à jour du fichier via Shiny
$wb = $excel.ActiveWorkbook

# --- "pro" sheet: update B2:B26 with new production values ---
$wsPro = $wb.Worksheets.Item("pro")
$proValues = @(2858042.4441151237,3234474.5079740775,3565093.4883621228,3529880.4330079281,3419058.0036565899,3483015.9777827743,3447876.9166528936,3260423.8878553719,3437118.5208715466,3592483.8895612396,3664033.8556317249,3770294.4123181165,4019364.5799250421,4142978.9348825701,4498682.490473724,4381460.758833277,5131718.9554420924,5378021.9417955615,5722553.4582232181,6097117,6710446,6780085.1173277302,7567508.2198896352,7885018.8946692916,7712467.9359351657)
for ($i = 0; $i -lt $proValues.Length; $i++) {
    $wsPro.Cells.Item($i + 2, 2).Value = $proValues[$i]
}

# --- "ind" sheet: update B2:B101 with new index values ---
$wsInd = $wb.Worksheets.Item("ind")
$indValues = @(936105.75982185255,928314.37579718232,932669.06303654052,937023.75027589116,1012244.4912122109,1029987.4107759548,1066009.9648907091,1117690.1755524143,1122695.6884051363,1164527.8223500301,1185286.8803847947,1185384.6609280694,1142475.8593310704,1138414.8840670069,1150661.2880954784,1180336.1719917036,1084249.3612124999,1114894.2239582008,1132674.5153484996,1135277.4777463179,1151487.9784192315,1139603.8948753451,1131106.8244867884,1128459.7750324016,1124606.6865500417,1127104.6963498918,1127520.6918648488,1125516.2250042984,1060659.6030455015,1061055.2880089344,1065499.9892827943,1072620.7577635364,1135770.117537596,1134113.0043248208,1121631.3875612176,1099177.6432509285,1205597.8928211923,1177701.6191501068,1158954.1335180823,1151427.737185688,1186814.5152755342,1193422.6170181904,1200653.74433066,1206272.5518997707,1228429.9151612595,1230550.1387321956,1232325.0163338482,1234690.7683964064,1300344.1974023965,1306926.6112458303,1316261.5987988391,1327880.6201230814,1337795.8748758314,1349740.085446263,1359463.9926887949,1365918.7108632801,1474849.0195391856,1473138.895730603,1467989.4191516854,1461678.0430899644,1417724.470969155,1420877.1067407792,1432915.5349530273,1452984.7709540613,1644827.0404468474,1669898.0884452055,1688890.5920948051,1701119.6272407076,1737779.2511036643,1747416.3848487742,1761462.435642438,1779879.0653348616,1736877.0942185703,1961089.3789478329,1896613.4740806296,1882097.3194965683,2121171.2827441031,2495512.1187999076,2455316.9204027979,2447788.9241538718,3301722.7553905826,3566202.869179457,3438202.2537474823,3550332.6303950218,4195514.327750179,4615164.8825397799,4354898.6114062462,4567212.970366247,4566786.9895859035,4963362.9757220447,4657293.9943095064,4942615.2227526931,4712152.7337034941,5322462.2117269132,4960914.7125555109,5283660.3261204222,4868591.7506114114,5221195.6109139081,5332671.5693537202,5241119.0438278206)
for ($i = 0; $i -lt $indValues.Length; $i++) {
    $wsInd.Cells.Item($i + 2, 2).Value = $indValues[$i]
}

# --- "conso" sheet: update B2:B26 with new consumption values ---
$wsConso = $wb.Worksheets.Item("conso")
$consoValues = @(295816.99654894491,370692.36013302696,411280.59703711804,392216.01327696169,363139.46666770708,372227.31880729255,398809.06466458761,363031.42136684904,367111.82403399161,400407.17527251417,411706.45971792488,427712.69536105386,456188.29991258855,472760.6916413203,517838.6058062659,505303.32377236773,604630.74908420851,644390.54540660861,694716.70748508919,734777,848205,885602.63662562764,961786.99276901514,1002140.7827081879,980210.54320555739)
for ($i = 0; $i -lt $consoValues.Length; $i++) {
    $wsConso.Cells.Item($i + 2, 2).Value = $consoValues[$i]
}

# "VA" sheet (VA = pro - conso) recalculates automatically via its formulas.

# --- Update selections / view state on every sheet to match the saved workbook ---
$wsVA = $wb.Worksheets.Item("VA")

$wsInd.Activate()
$excel.ActiveWindow.ScrollRow = 96
$wsInd.Range("B102:B113").Select()

$wsVA.Activate()
$wsVA.Range("B102:B113").Select()

$wsConso.Activate()
$wsConso.Range("B102:B113").Select()

# "pro" is the sheet that was active/selected when the workbook was saved
$wsPro.Activate()
$wsPro.Range("B102:B113").Select()

